$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.2526690391459075
$ws1.Range("C2").Value = 0.0625
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.1176470588235294
$ws1.Range("F2").Value = 0.25
$ws1.Range("G2").Value = 0.6341463414634146
$ws1.Range("H2").Value = 0.77919341894061
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 420
$ws1.Range("K2").Value = 114
$ws1.Range("L2").Value = 0

# ---- Sheet 2: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 -> label "0"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.2134831460674157
$ws2.Range("D2").Value = 0.3518518518518519

# row 3 -> label "1"
$ws2.Range("B3").Value = 0.0625
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.1176470588235294

# row 4 -> label "accuracy"
$ws2.Range("B4").Value = 0.2526690391459075
$ws2.Range("C4").Value = 0.2526690391459075
$ws2.Range("D4").Value = 0.2526690391459075
$ws2.Range("E4").Value = 0.2526690391459075

# row 5 -> label "macro avg"
$ws2.Range("B5").Value = 0.53125
$ws2.Range("C5").Value = 0.6067415730337079
$ws2.Range("D5").Value = 0.2347494553376906

# row 6 -> label "weighted avg"
$ws2.Range("B6").Value = 0.9532918149466192
$ws2.Range("C6").Value = 0.2526690391459075
$ws2.Range("D6").Value = 0.3401832856511525

# ---- Sheet 3: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 114
$ws3.Range("C2").Value = 420
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
